$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14 ("Percentile Lingo") - Content Placeholder 2
# Merge the 3 trailing runs of the last paragraph into a single run:
#   ' percentile” is ' + 'more obvious ' + 'as “my percentile rank is 92” '
#   -> ' percentile” is more obvious as “my percentile rank is 92” '
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

# The last paragraph in this placeholder is paragraph 8.
$lastPara = $tr14.Paragraphs(8, 1)

# Locate the run boundary right after the superscript "nd" run: the
# merged text itself begins with ' percentile...' so we can find it
# reliably inside the paragraph text.
$paraText = $lastPara.Text
$markerIdx = $paraText.IndexOf(" percentile")
$mergeStart = $lastPara.Start + $markerIdx
$mergeLength = $lastPara.Start + $lastPara.Length - $mergeStart

$mergeRange = $tr14.Characters($mergeStart, $mergeLength)
$mergeRange.Text = " percentile” is more obvious as “my percentile rank is 92” "

# ---------------------------------------------------------------------------
# Slide 16 ("IQR and Ranges") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)
$tr16 = $sh16.TextFrame.TextRange

# 1) "Can also help with automated outlier filter. " ->
#    "Can also help with automated outlier ident./filter. "
$paraOutlier = $tr16.Paragraphs(6, 1)
$cOutlier = $tr16.Characters($paraOutlier.Start, $paraOutlier.Length)
$cOutlier.Text = "Can also help with automated outlier ident./filter. "

# 2) Insert two new level-1 (sub-bullet) paragraphs right after it.
$tr16b = $sh16.TextFrame.TextRange
$paraOutlier2 = $tr16b.Paragraphs(6, 1)
$inserted = $paraOutlier2.InsertAfter("`rN multiple of IQR as a cutoff. `rSmart filtering based on the actual data. ")

$tr16c = $sh16.TextFrame.TextRange
$newPara1 = $tr16c.Paragraphs(7, 1)
$newPara1.IndentLevel = 2
$newPara2 = $tr16c.Paragraphs(8, 1)
$newPara2.IndentLevel = 2

# 3) "What if the data is not normal?" -> "What if the data is not normal…"
$tr16d = $sh16.TextFrame.TextRange
$paraNormal = $tr16d.Paragraphs(9, 1)
$cNormal = $tr16d.Characters($paraNormal.Start, $paraNormal.Length)
$cNormal.Text = "What if the data is not normal…"
